# Aula9 - Arquitetura.pptx edit
# Commit: "PPTs Arq. de Comp. Yduqs 2022_2 - 27102022"
#
# Replaces the outgoing professor's name ("Talita Rocha Pinheiro") with the
# incoming one ("Heleno Cardoso") on the title slide and on the closing
# slide, and adds a credits line on the closing slide crediting the
# previous professor.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 1 (title slide): "Talita Rocha Pinheiro" -> "Heleno Cardoso"
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$tf1 = $s1.Shapes.Item(1).TextFrame
$tr1 = $tf1.TextRange
$nameRun1 = $tr1.Paragraphs(2).Runs(1)
$nameRun1.Text = "Heleno Cardoso"

# ---------------------------------------------------------------------
# Last slide (closing slide): same name swap, plus a new credits line
# thanking the previous professor.
# ---------------------------------------------------------------------
$lastIndex = $p.Slides.Count
$sLast = $p.Slides.Item($lastIndex)
$tfLast = $sLast.Shapes.Item(1).TextFrame
$trLast = $tfLast.TextRange
$nameRunLast = $trLast.Paragraphs(2).Runs(1)
$nameRunLast.Text = "Heleno Cardoso"

# Append a brand-new paragraph after the name, carrying the same
# bold/size/color/font formatting, crediting the previous professor.
$null = $trLast.InsertAfter("`rCréditos: Professora Talita Rocha Pinheiro")
